# Updates 展览 (sheet1) and 全部类型 (sheet4) with the latest event listings.
# Row 2 (南昌·DSL国风动漫游戏嘉年华) moves to the bottom (new row 17) with an updated
# date/time window; rows 3-14 shift up by one; a new row (抚州·第七届FZ动漫文化节) is
# inserted at row 14; rows 15-16 are re-affirmed unchanged.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
foreach ($sheetName in $sheetNames) {
  $ws = $wb.Worksheets.Item($sheetName)

  # Give the new row 17 the same style as the existing index column (bold/border/center).
  $ws.Range("A16").Copy()
  $ws.Range("A17").PasteSpecial(-4122)

  # Column B holds plain "yyyy-mm-dd" text, not real dates; force text format first so
  # Excel does not silently coerce the assigned strings into date serials.
  $ws.Range("B2:B17").NumberFormat = "@"

  # Row 2: 景德镇·江报国风动漫展 
  $ws.Range("A2").Value = 1
  $ws.Range("B2").Value = "2024-03-09"
  $ws.Range("C2").Value = "景德镇·江报国风动漫展 "
  $ws.Range("D2").Value = "迎宾大道与寺山路交叉口东200米 陶博城"
  $ws.Range("E2").Value = "2024.03.09 09:00-03.10 17:00"
  $ws.Range("F2").Value = 1062
  $ws.Range("G2").Value = 55
  $ws.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=81362"
  $ws.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202402/hdlmhoLp1708932790894.jpeg"

  # Row 3: 景德镇·原神X崩铁X崩坏动漫展only
  $ws.Range("A3").Value = 2
  $ws.Range("B3").Value = "2024-03-16"
  $ws.Range("C3").Value = "景德镇·原神X崩铁X崩坏动漫展only"
  $ws.Range("D3").Value = "陶阳南路188号 晨枫臻品酒店"
  $ws.Range("E3").Value = "2024.03.16 10:00-03.16 17:00"
  $ws.Range("F3").Value = 74
  $ws.Range("G3").Value = 55
  $ws.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=80920"
  $ws.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202401/IugBckTp1705469476482.png"

  # Row 4: 江西·ShiningStaR动漫游戏文化节5th
  $ws.Range("A4").Value = 3
  $ws.Range("B4").Value = "2024-03-16"
  $ws.Range("C4").Value = "江西·ShiningStaR动漫游戏文化节5th"
  $ws.Range("D4").Value = "江西科技学院内 江西科技学院体育馆"
  $ws.Range("E4").Value = "2024.03.16 09:30-03.17 17:00"
  $ws.Range("F4").Value = 3055
  $ws.Range("G4").Value = "不可售"
  $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=81792"
  $ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202403/p3TpZeAQ1709544877660.jpeg"

  # Row 5: 上饶·原×铁×崩only
  $ws.Range("A5").Value = 4
  $ws.Range("B5").Value = "2024-03-23"
  $ws.Range("C5").Value = "上饶·原×铁×崩only"
  $ws.Range("D5").Value = "五三东大道42号 回禾酒店"
  $ws.Range("E5").Value = "2024.03.23 10:00-03.23 17:00"
  $ws.Range("F5").Value = 41
  $ws.Range("G5").Value = 60
  $ws.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=81103"
  $ws.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202401/pp6c5TsC1705647180602.jpeg"

  # Row 6: 南昌·AP动漫游戏嘉年华
  $ws.Range("A6").Value = 5
  $ws.Range("B6").Value = "2024-03-23"
  $ws.Range("C6").Value = "南昌·AP动漫游戏嘉年华"
  $ws.Range("D6").Value = "八一桥街道青山南路118号 蓝海会展中心"
  $ws.Range("E6").Value = "2024.03.23 09:00-03.24 17:00"
  $ws.Range("F6").Value = 2218
  $ws.Range("G6").Value = 58.5
  $ws.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=81232"
  $ws.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202401/NZv97SmS1705912230957.jpeg"

  # Row 7: 南昌·运动番only春季集训（取消）
  $ws.Range("A7").Value = 6
  $ws.Range("B7").Value = "2024-03-23"
  $ws.Range("C7").Value = "南昌·运动番only春季集训（取消）"
  $ws.Range("D7").Value = "创新三路777号 南昌小飞侠章鱼文化体育公园"
  $ws.Range("E7").Value = "2024.03.23 10:00-03.24 17:00"
  $ws.Range("F7").Value = 188
  $ws.Range("G7").Value = "不可售"
  $ws.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=81950"
  $ws.Range("I7").Value = "//i1.hdslb.com/bfs/openplatform/202402/bm4uH4qB1708425538357.jpeg"

  # Row 8: 南昌·AP动漫游戏  嘉年华内场票-小N&子音
  $ws.Range("A8").Value = 7
  $ws.Range("B8").Value = "2024-03-24"
  $ws.Range("C8").Value = "南昌·AP动漫游戏  嘉年华内场票-小N&子音"
  $ws.Range("D8").Value = "八一桥街道青山南路118号 蓝海会展中心"
  $ws.Range("E8").Value = "2024.03.24 09:00-03.24 17:00"
  $ws.Range("F8").Value = 110
  $ws.Range("G8").Value = 218
  $ws.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=81973"
  $ws.Range("I8").Value = "//i0.hdslb.com/bfs/openplatform/202402/zbG5HICL1708504962467.jpeg"

  # Row 9: 南昌·CM01动漫游戏博览会
  $ws.Range("A9").Value = 8
  $ws.Range("B9").Value = "2024-03-30"
  $ws.Range("C9").Value = "南昌·CM01动漫游戏博览会"
  $ws.Range("D9").Value = "怀玉山大道1315号 南昌绿地国际博览中心"
  $ws.Range("E9").Value = "2024.03.30 10:00-03.31 17:00"
  $ws.Range("F9").Value = 1048
  $ws.Range("G9").Value = 55
  $ws.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=81691"
  $ws.Range("I9").Value = "//i2.hdslb.com/bfs/openplatform/202402/9cMJMElF1708938074308.png"

  # Row 10: 鹰潭·原×铁×崩only
  $ws.Range("A10").Value = 9
  $ws.Range("B10").Value = "2024-03-30"
  $ws.Range("C10").Value = "鹰潭·原×铁×崩only"
  $ws.Range("D10").Value = "南站路24号 回禾酒店(鹰潭火车站店)"
  $ws.Range("E10").Value = "2024.03.30 10:00-03.30 17:00"
  $ws.Range("F10").Value = 33
  $ws.Range("G10").Value = 60
  $ws.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=81097"
  $ws.Range("I10").Value = "//i2.hdslb.com/bfs/openplatform/202401/q0AZaXAk1705646244207.jpeg"

  # Row 11: 新余·文旅国漫嘉年华暨BM次元盛典
  $ws.Range("A11").Value = 10
  $ws.Range("B11").Value = "2024-03-31"
  $ws.Range("C11").Value = "新余·文旅国漫嘉年华暨BM次元盛典"
  $ws.Range("D11").Value = "五一南路与仙女湖大道交叉口西北 老上海风情街白金汉宫"
  $ws.Range("E11").Value = "2024.03.31 10:00-03.31 17:00"
  $ws.Range("F11").Value = 37
  $ws.Range("G11").Value = 60
  $ws.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=82208"
  $ws.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202403/aXc6vPDP1709547191851.jpeg"

  # Row 12: 赣州·第三届半夏动漫展
  $ws.Range("A12").Value = 11
  $ws.Range("B12").Value = "2024-04-04"
  $ws.Range("C12").Value = "赣州·第三届半夏动漫展"
  $ws.Range("D12").Value = "105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心"
  $ws.Range("E12").Value = "2024.04.04 10:00-04.06 17:00"
  $ws.Range("F12").Value = 256
  $ws.Range("G12").Value = 50
  $ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=82235"
  $ws.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202403/4DWZWYGm1709278879159.jpeg"

  # Row 13: 赣州·赣次元·归来国风动漫节
  $ws.Range("A13").Value = 12
  $ws.Range("B13").Value = "2024-04-04"
  $ws.Range("C13").Value = "赣州·赣次元·归来国风动漫节"
  $ws.Range("D13").Value = "客家大道568号文清外国语学校旁 赣州市文清外国语学校国际交流中心"
  $ws.Range("E13").Value = "2024.04.04 10:00-04.04 17:00"
  $ws.Range("F13").Value = 93
  $ws.Range("G13").Value = 40
  $ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=82125"
  $ws.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202402/8RNepTak1709022774421.jpeg"

  # Row 14: 抚州·第七届FZ动漫文化节
  $ws.Range("A14").Value = 13
  $ws.Range("B14").Value = "2024-04-05"
  $ws.Range("C14").Value = "抚州·第七届FZ动漫文化节"
  $ws.Range("D14").Value = "迎宾大道288号 凤凰世纪名都大酒店"
  $ws.Range("E14").Value = "2024.04.05 09:30-04.05 17:00"
  $ws.Range("F14").Value = 2
  $ws.Range("G14").Value = 50
  $ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=82381"
  $ws.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202403/Y725SN0G1709694367526.jpeg"

  # Row 15: 南昌·原X穹X崩only
  $ws.Range("A15").Value = 14
  $ws.Range("B15").Value = "2024-04-13"
  $ws.Range("C15").Value = "南昌·原X穹X崩only"
  $ws.Range("D15").Value = "丰和北大道299号 新吉花园酒店"
  $ws.Range("E15").Value = "2024.04.13 10:00-04.13 17:00"
  $ws.Range("F15").Value = 96
  $ws.Range("G15").Value = 65
  $ws.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=80807"
  $ws.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202402/kfK13XvH1709202705153.jpeg"

  # Row 16: 南昌·第二届漫拥动漫嘉年华mini
  $ws.Range("A16").Value = 15
  $ws.Range("B16").Value = "2024-04-13"
  $ws.Range("C16").Value = "南昌·第二届漫拥动漫嘉年华mini"
  $ws.Range("D16").Value = "南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆"
  $ws.Range("E16").Value = "2024.04.13 10:00-04.14 18:00"
  $ws.Range("F16").Value = 44
  $ws.Range("G16").Value = 39.9
  $ws.Range("H16").Value = "https://show.bilibili.com/platform/detail.html?id=82210"
  $ws.Range("I16").Value = "//i0.hdslb.com/bfs/openplatform/202402/KYd0bfk11709203777701.png"

  # Row 17: 南昌·DSL国风动漫游戏嘉年华
  $ws.Range("A17").Value = 16
  $ws.Range("B17").Value = "2024-04-20"
  $ws.Range("C17").Value = "南昌·DSL国风动漫游戏嘉年华"
  $ws.Range("D17").Value = "沿江北路69号 瑞颐大酒店"
  $ws.Range("E17").Value = "2024.04.20 09:00-04.21 17:00"
  $ws.Range("F17").Value = 73
  $ws.Range("G17").Value = 35
  $ws.Range("H17").Value = "https://show.bilibili.com/platform/detail.html?id=82107"
  $ws.Range("I17").Value = "//i0.hdslb.com/bfs/openplatform/202402/QDlumVb41708943318282.jpeg"

}
